# Updates cryptos list data (prices / volume / rows) per latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: writes a value and forces it to stay a text cell (matches the
# workbook convention of storing all Price/Volume figures as plain text),
# then clears the scratch number-format back off the cell style so no
# stray formatting is left behind.
function Set-TextValue($addr, $value) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "64.016.38"
Set-TextValue "E2" "  +1.55%  "
# Row 3
Set-TextValue "D3" "3.228.19"
Set-TextValue "E3" "  -1.42%  "
# Row 4
Set-TextValue "E4" "  +0.02%  "
# Row 5
Set-TextValue "D5" "595.44"
Set-TextValue "E5" "  -0.52%  "
# Row 6
Set-TextValue "D6" "137.97"
Set-TextValue "E6" "  -0.29%  "
# Row 7
Set-TextValue "E7" "  -0.10%  "
# Row 8
Set-TextValue "D8" "3.226.55"
Set-TextValue "E8" "  -1.47%  "
# Row 9
Set-TextValue "E9" "  +1.10%  "
# Row 10
Set-TextValue "E10" "  -1.90%  "
# Row 11
Set-TextValue "D11" "5.34"
Set-TextValue "E11" "  -2.08%  "
# Row 12
Set-TextValue "E12" "  -0.89%  "
# Row 13
Set-TextValue "E13" "  +0.09%  "
# Row 14
Set-TextValue "D14" "35.20"
Set-TextValue "E14" "  +3.58%  "
# Row 15
Set-TextValue "D15" "3.762.26"
Set-TextValue "E15" "  -1.45%  "
# Row 16
Set-TextValue "E16" "  -1.79%  "
# Row 17
Set-TextValue "D17" "3.222.03"
Set-TextValue "E17" "  -1.71%  "
# Row 18
Set-TextValue "D18" "64.065.21"
Set-TextValue "E18" "  +1.45%  "
# Row 19
Set-TextValue "D19" "6.64"
Set-TextValue "E19" "  -1.46%  "
# Row 20
Set-TextValue "D20" "468.27"
Set-TextValue "E20" "  -0.80%  "
# Row 21
Set-TextValue "D21" "14.13"
Set-TextValue "E21" "  +2.40%  "
# Row 22
Set-TextValue "D22" "0.709"
Set-TextValue "E22" "  -2.25%  "
# Row 23
Set-TextValue "D23" "7.75"
Set-TextValue "E23" "  -1.33%  "
# Row 24
Set-TextValue "D24" "13.52"
Set-TextValue "E24" "  -0.74%  "
# Row 25
Set-TextValue "D25" "83.62"
Set-TextValue "E25" "  -0.76%  "
# Row 26
Set-TextValue "E26" "  +0.26%  "
# Row 27
Set-TextValue "E27" "  -1.21%  "
# Row 28
Set-TextValue "D28" "0.997"
Set-TextValue "E28" "  -0.38%  "
# Row 29
Set-TextValue "D29" "7.89"
Set-TextValue "E29" "  -1.01%  "
# Row 30
Set-TextValue "B30" "ImmutableX"
Set-TextValue "C30" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D30" "2.10"
Set-TextValue "E30" "  -0.54%  "
# Row 31
Set-TextValue "B31" "NEARProtocol"
Set-TextValue "C31" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D31" "6.88"
Set-TextValue "E31" "  -2.50%  "
# Row 32
Set-TextValue "D32" "27.71"
Set-TextValue "E32" "  -1.49%  "
# Row 33
Set-TextValue "E33" "  -0.08%  "
# Row 34
Set-TextValue "D34" "2.45"
Set-TextValue "E34" "  -0.41%  "
# Row 35
Set-TextValue "D35" "1.05"
Set-TextValue "E35" "  -2.99%  "
# Row 36
Set-TextValue "D36" "5.95"
Set-TextValue "E36" "  +0.47%  "
# Row 37
Set-TextValue "D37" "51.73"
Set-TextValue "E37" "  -0.22%  "
# Row 38
Set-TextValue "D38" "0.0₃0735"
Set-TextValue "E38" "  +2.16%  "
# Row 39
Set-TextValue "D39" "0.0396"
Set-TextValue "E39" "  +0.86%  "
# Row 40
Set-TextValue "D40" "2.79"
Set-TextValue "E40" "  +5.46%  "
# Row 41
Set-TextValue "D41" "407.42"
Set-TextValue "E41" "  -3.18%  "
# Row 42
Set-TextValue "D42" "8.17"
Set-TextValue "E42" "  -0.30%  "
# Row 43
Set-TextValue "D43" "0.114"
Set-TextValue "E43" "  -2.68%  "
# Row 44
Set-TextValue "D44" "2.854.19"
Set-TextValue "E44" "  -7.45%  "
# Row 45
Set-TextValue "E45" "  -0.45%  "
# Row 46
Set-TextValue "D46" "2.17"
Set-TextValue "E46" "  +0.12%  "
# Row 47
Set-TextValue "D47" "36.31"
Set-TextValue "E47" "  +1.76%  "
# Row 48
Set-TextValue "B48" "Monero"
Set-TextValue "C48" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D48" "127.63"
Set-TextValue "E48" "  +0.62%  "
# Row 49
Set-TextValue "B49" "USDe"
Set-TextValue "C49" "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue "D49" "0.998"
Set-TextValue "E49" "  -0.04%  "
# Row 50
Set-TextValue "D50" "25.88"
Set-TextValue "E50" "  +0.21%  "
# Row 51
Set-TextValue "D51" "0.113"
Set-TextValue "E51" "  +0.22%  "
